$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d40 = @'
<br><a style="color: red; font-weight: bold;">Volvo Internal Use Only</a>
<br>
<br>차량은 간혹 문제가 발생하거나 발생하기 전 볼보 자동차에서 미리 해결하기 위한 문제 해결책으로 리콜이나 무상수리를 진행할 수 있습니다.
<br>이 문서에서는 리콜 및 무상 수리에 대한 정보를 확인할 수 있습니다.
<br>
<br>볼보자동차코리아는 대한민국의 법률을 준수하며 차량을 제조하고 대한민국에서 판매될 때 대한민국 법률을 준수합니다.
<br>
<br>리콜이나 무상수리나 리콜의 대상인 경우 KATRI(자동차 안전 연구원)을 통해 고객에게 메시지나 우편으로 고객에게 통보됩니다.
<br>
<br>무상 수리 및 리콜을 Quality Bulletin(QB)으로 표현합니다.
<br>문제가 있는 경우 스웨덴 볼보 본사에서 조사를 진행하게 됩니다.
<br>
<br><h3>QB의 종류</h3>
<br><strong>・ QB-D (Delivery Stop) - 딜리버리 중지</strong>
<br>    - PDI센터를 통해 입고 후 조치를 진행 고객에게 출고를 진행하지 않도록 중지하는 것입니다.
<br>
<br><strong>・ QB-S (Service Action) - 서비스 조치</strong>
<br>    - QB-S/QB-S=한국에서는 동일한 무상수리로 간주합니다. 이름만 다르며, 고객 만족 차원에서 진행합니다. 반드시 첫 서비스센터 방문 시 수행되어야 하는 항목입니다.
<br>    - 고객에게 10일 이내 통지 필요한 항목입니다.
<br>
<br><strong>・ QB-A (Service Campaigns) - 서비스 캠페인</strong>
<br>    - 전사적으로 적극적으로 고객을 불러들여야 하는 QB-S에 비해 심각한 조치입니다. QB-S보다 좀더 적극적인 입고 유도가 필요한 항목입니다.
<br>    - 고객에게 10일 이내 통지 필요한 항목입니다.
<br>
<br><strong>・ QB-R (Recall) 리콜</strong>
<br>    - 법규를 위반한 경우 / 안전 사양에 문제가 있는 경우에 진행 합니다.
<br>    - 리콜 조치는 다른 조치들보다 최우선으로 진행 됩니다.
<br>    - 리콜 개시일로부터 30일 이내 고객에게 통지가 완료되어야 합니다.
<br>    - 리콜의 경우 사진이나 이미지를 통지문에 첨부해야 합니다.
<br>
<br><em>대체로 국토교통부 자동차 정책과에서 문의가 인입될 수 있습니다. 이를 포함한 정부 기관에서 문의하는 경우 적절한 곳으로 에스컬레이션이 진행되어야 합니다.</em>
<br>리콜에 대해서 VCK에 문의할 때 R로 시작하는 번호로 소통 필요하며, 무상 수리의 경우 S나 A로 시작하는 번호로 소통 필요합니다.
'@

$ws.Range("D40").Value = $d40

$d38 = @'
고객이 Hej Volvo 앱의 차량 정보 동기화를 위해 문의하는 경우 아래 기준으로 진행합니다.
<br>
<br><string style="color: red;"><em>고객지원팀에서는 절대 실운전자 정보를 수정하면 안 됩니다.</em></string>
<br><string style="color: red;"><em>고객지원팀은 계약자 정보만 수정하며, 정보 변경이 필요한 경우 서비스센터를 방문했으나, 로그인이 안되는 경우로, 전화번호나 이름 철자 하나가 틀린 경우에만 수정이 가능합니다. 나머지는 서비스센터로 안내합니다.</em></string>
<br><h3>계약자 정보 변경을 위해 다음 정보를 고객으로부터 수집한 후 수정합니다.</h3>
<br>1. 실운전자 정보 (실 운전자에게 Callback하여 해당 내용을 확인해야 합니다.)
<br>2. 차량등록증 (차량등록증 상에 리스사로 되어있는 경우 차량등록증+차대번호와 실운전자가 기재된 리스계약서)
<br>
<br>고객으로부터 서류를 수집하는 경우 반드시 주민번호 전체를 가리고 보내줄 수 있도록 안내합니다.
<br>계약자 정보를 수집한 후에는 반드시 해당 서류를 파기할 수 있도록 합니다.
'@

$ws.Range("D38").Value = $d38

$d41 = @'
<a style="color: red; font-weight: bold;">Volvo Internal Use Only</a>
<br>
<br>이곳에 내용이 들어갈 예정 :)
'@

$ws.Range("D41").Value = $d41

$ws.Range("C51").Select()
